# Restore the original input data for validation.
# The sheet previously had row 2 with a stray B2 value of 2 and was
# missing rows 3-5 of the ID_electricity / id_electricity_feed_in /
# id_gases / price_unit pattern. Put the data back the way it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 2 (B2 was accidentally 2, should be 1)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "cent/Wh"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "cent/Wh"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "cent/Wh"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "cent/Wh"

# Re-fit the columns that hold the new, wider data and leave the
# selection where the user ended up after typing the data in.
$ws.Columns.Item(1).ColumnWidth = 13.45
$ws.Columns.Item(2).ColumnWidth = 11.7
$ws.Columns.Item(3).ColumnWidth = 19.7

$ws.Range("A6").Select()
